# "Generate Report for Handoff" — refresh the handoff timestamps for the
# 3389134b-1eb0-4cc6-8958-a0fe99b63eb7 record (row 6 on every sheet) to
# reflect a freshly (re-)generated handoff xliff.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G), row 6 ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G6").Value = "2016-08-12 14:49:26"

# --- zh-cn sheet: "Latest Handoff Datetime" column (H), row 6 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H6").Value = "2016-08-12 14:49:18"

# --- de-de sheet: "Latest Handoff Datetime" column (H), row 6 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H6").Value = "2016-08-12 14:49:26"
